$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 27784000
$ws.Range("I116").Value = 62501748
$ws.Range("K116").Value = 62501748
$ws.Range("M116").Value = -62498306
$ws.Range("H123").Value = 50228.285
$ws.Range("J123").Value = 50228.285
$ws.Range("L123").Value = 50228.285
$ws.Range("N123").Value = -60028.285
$ws.Range("H131").Value = 4065.4443
$ws.Range("I131").Value = 3431.8333
$ws.Range("K131").Value = 10295.4999
$ws.Range("M131").Value = -5255.499899999999
$ws.Range("H132").Value = 1532.3903
$ws.Range("I132").Value = 1464.4865
$ws.Range("K132").Value = 4393.4595
$ws.Range("M132").Value = -1863.4595
$ws.Range("H137").Value = 2758.7437
$ws.Range("J137").Value = 2899.6667
$ws.Range("L137").Value = 8699.000100000001
$ws.Range("N137").Value = -13799.0001
$ws.Range("H138").Value = 5161.3213
$ws.Range("J138").Value = 7426.914
$ws.Range("L138").Value = 22280.742
$ws.Range("N138").Value = -32560.742
$ws.Range("H141").Value = 5748438
$ws.Range("I141").Value = 6945226
$ws.Range("J141").Value = 3856
$ws.Range("K141").Value = 20835678
$ws.Range("L141").Value = 11568
$ws.Range("M141").Value = -20830498
$ws.Range("N141").Value = -21928
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2057216.1
$ws.Range("I32").Value = 2279414.2
$ws.Range("J32").Value = 20399.334
$ws.Range("K32").Value = 2279414.2
$ws.Range("L32").Value = 20399.334
$ws.Range("M32").Value = -2279127.2
$ws.Range("N32").Value = -20973.334
$ws.Range("H45").Value = 4392.36
$ws.Range("I45").Value = 3594.5264
$ws.Range("J45").Value = 6918.8335
$ws.Range("K45").Value = 3594.5264
$ws.Range("L45").Value = 6918.8335
$ws.Range("M45").Value = -3217.5264
$ws.Range("N45").Value = -7672.8335
$ws.Range("H74").Value = 31268.676
$ws.Range("I74").Value = 41104.88
$ws.Range("J74").Value = 3945.889
$ws.Range("K74").Value = 41104.88
$ws.Range("L74").Value = 3945.889
$ws.Range("M74").Value = -40230.88
$ws.Range("N74").Value = -5693.889
$ws.Range("H77").Value = 31268.676
$ws.Range("I77").Value = 41104.88
$ws.Range("J77").Value = 3945.889
$ws.Range("K77").Value = 205524.4
$ws.Range("L77").Value = 19729.445
$ws.Range("M77").Value = -201156.4
$ws.Range("N77").Value = -28465.445
$ws.Range("H110").Value = 47620564
$ws.Range("I110").Value = 1766.5
$ws.Range("K110").Value = 1766.5
$ws.Range("M110").Value = 278.5
$ws.Range("H122").Value = 3107.72
$ws.Range("I122").Value = 1794.9524
$ws.Range("K122").Value = 5384.857199999999
$ws.Range("M122").Value = -2934.857199999999
$ws.Range("H132").Value = 4477.2686
$ws.Range("I132").Value = 3261.7073
$ws.Range("J132").Value = 6394.115
$ws.Range("K132").Value = 9785.1219
$ws.Range("L132").Value = 19182.345
$ws.Range("M132").Value = -7255.1219
$ws.Range("N132").Value = -24242.345
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 29999
$ws.Range("J76").Value = 29999
$ws.Range("L76").Value = 29999
$ws.Range("N76").Value = -30629
$ws.Range("H79").Value = 29999
$ws.Range("J79").Value = 29999
$ws.Range("L79").Value = 29999
$ws.Range("N79").Value = -32183
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H105").Value = 1054366
$ws.Range("I105").Value = 1112692
$ws.Range("K105").Value = 1112692
$ws.Range("M105").Value = -1110945
$ws.Range("H132").Value = 65000
$ws.Range("J132").Value = 65000
$ws.Range("L132").Value = 65000
$ws.Range("N132").Value = -75120
$ws.Range("H134").Value = 5982.5435
$ws.Range("I134").Value = 2156.9546
$ws.Range("K134").Value = 6470.8638
$ws.Range("M134").Value = -3935.8638
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6279.2812
$ws.Range("I31").Value = 2833.7368
$ws.Range("J31").Value = 11315.077
$ws.Range("K31").Value = 2833.7368
$ws.Range("L31").Value = 11315.077
$ws.Range("M31").Value = -2538.7368
$ws.Range("N31").Value = -11905.077
$ws.Range("H34").Value = 6279.2812
$ws.Range("I34").Value = 2833.7368
$ws.Range("J34").Value = 11315.077
$ws.Range("K34").Value = 2833.7368
$ws.Range("L34").Value = 11315.077
$ws.Range("M34").Value = -2631.7368
$ws.Range("N34").Value = -11719.077
$ws.Range("H62").Value = 7753.857
$ws.Range("I62").Value = 6266
$ws.Range("K62").Value = 6266
$ws.Range("M62").Value = -5642
$ws.Range("H65").Value = 7753.857
$ws.Range("I65").Value = 6266
$ws.Range("K65").Value = 31330
$ws.Range("M65").Value = -28210
$ws.Range("H99").Value = 6766.4
$ws.Range("I99").Value = 4997.1665
$ws.Range("K99").Value = 4997.1665
$ws.Range("M99").Value = -3499.1665
$ws.Range("H126").Value = 6766.4
$ws.Range("I126").Value = 4997.1665
$ws.Range("K126").Value = 14991.4995
$ws.Range("M126").Value = -12521.4995
$ws.Range("H134").Value = 4270.9844
$ws.Range("I134").Value = 1792.8223
$ws.Range("K134").Value = 5378.4669
$ws.Range("M134").Value = -2843.4669
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 208.13792
$ws.Range("I23").Value = 153.72728
$ws.Range("K23").Value = 461.18184
$ws.Range("M23").Value = -226.18184
$ws.Range("H55").Value = 68259960
$ws.Range("J55").Value = 5889606.5
$ws.Range("L55").Value = 17668819.5
$ws.Range("N55").Value = -17669173.5
$ws.Range("H136").Value = 3441.7778
$ws.Range("I136").Value = 1130.1333
$ws.Range("K136").Value = 3390.3999
$ws.Range("M136").Value = 1709.6001
$ws.Range("H138").Value = 46692.625
$ws.Range("I138").Value = 54232.6
$ws.Range("K138").Value = 162697.8
$ws.Range("M138").Value = -157557.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3246
$ws.Range("I80").Value = 3246
$ws.Range("K80").Value = 3246
$ws.Range("M80").Value = -2248
$ws.Range("H83").Value = 3246
$ws.Range("I83").Value = 3246
$ws.Range("K83").Value = 16230
$ws.Range("M83").Value = -11238
$ws.Range("H122").Value = 1907724.4
$ws.Range("I122").Value = 2898660
$ws.Range("J122").Value = 2079.077
$ws.Range("K122").Value = 8695980
$ws.Range("L122").Value = 6237.231000000001
$ws.Range("M122").Value = -8693530
$ws.Range("N122").Value = -11137.231
$ws.Range("H132").Value = 4745.7
$ws.Range("I132").Value = 2196.275
$ws.Range("J132").Value = 9844.549999999999
$ws.Range("K132").Value = 6588.825000000001
$ws.Range("L132").Value = 29533.65
$ws.Range("M132").Value = -4058.825000000001
$ws.Range("N132").Value = -34593.64999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 728.9
$ws.Range("J16").Value = 492
$ws.Range("L16").Value = 492
$ws.Range("N16").Value = -832
$ws.Range("H101").Value = 44263.6
$ws.Range("J101").Value = 44263.6
$ws.Range("L101").Value = 44263.6
$ws.Range("N101").Value = -50753.6
$ws.Range("H103").Value = 38397
$ws.Range("J103").Value = 38397
$ws.Range("L103").Value = 38397
$ws.Range("N103").Value = -40741
$ws.Range("H106").Value = 34813
$ws.Range("J106").Value = 34813
$ws.Range("L106").Value = 34813
$ws.Range("N106").Value = -37337
$ws.Range("H132").Value = 11118089
$ws.Range("I132").Value = 22730436
$ws.Range("J132").Value = 10628
$ws.Range("K132").Value = 68191308
$ws.Range("L132").Value = 31884
$ws.Range("M132").Value = -68188778
$ws.Range("N132").Value = -36944
$ws.Range("H136").Value = 8711.258
$ws.Range("I136").Value = 1542.5385
$ws.Range("K136").Value = 4627.6155
$ws.Range("M136").Value = -2077.6155
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 2619.4
$ws.Range("J41").Value = 2619.4
$ws.Range("L41").Value = 2619.4
$ws.Range("N41").Value = -3399.4
$ws.Range("H126").Value = 2793.7334
$ws.Range("J126").Value = 6166.6665
$ws.Range("L126").Value = 18499.9995
$ws.Range("N126").Value = -23439.9995
$ws.Range("H130").Value = 59379.5
$ws.Range("J130").Value = 59379.5
$ws.Range("L130").Value = 59379.5
$ws.Range("N130").Value = -69419.5
$ws.Range("H132").Value = 22747116
$ws.Range("I132").Value = 38471852
$ws.Range("K132").Value = 115415556
$ws.Range("M132").Value = -115413026
$ws.Range("H136").Value = 23282616
$ws.Range("I136").Value = 35715144
$ws.Range("J136").Value = 75229.60000000001
$ws.Range("K136").Value = 107145432
$ws.Range("L136").Value = 225688.8
$ws.Range("M136").Value = -107142882
$ws.Range("N136").Value = -230788.8
